$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------------
# New crime data collected - weekly CompStat report refresh (113th Precinct)
#   - Bulletin header: Volume/Number bumped, reporting week advanced by 7 days
#   - Crime-stat grid (rows 14-30, cols C:N) refreshed with new weekly figures
# ----------------------------------------------------------------------------

# --- Header text (rich-text shared strings -> plain text on write, value-equivalent) ---
$ws.Range("A8").Value = "Volume 32   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/10/2025  Through  2/16/2025"

# --- Column H got a touch narrower once Excel's bestFit recalculated for the new values ---
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# --- Step 1: fix cell styles for cells whose type changes (text<->number) ---
$ws.Range("C14").Copy($ws.Range("F14"))
$ws.Range("C16").Copy($ws.Range("C23"))
$ws.Range("C16").Copy($ws.Range("F23"))
$ws.Range("C16").Copy($ws.Range("I23"))
$ws.Range("K16").Copy($ws.Range("L23"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))

# --- Step 2: write numeric/text values ---
$ws.Range("N14").Value = -83.333333333333
$ws.Range("L15").Value = -66.666666666666
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("I16").Value = 11
$ws.Range("J16").Value = 16
$ws.Range("K16").Value = -31.25
$ws.Range("L16").Value = -21.428571428571
$ws.Range("M16").Value = -72.5
$ws.Range("N16").Value = -92.028985507246
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -28.571428571428
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 22.727272727272
$ws.Range("I17").Value = 46
$ws.Range("J17").Value = 45
$ws.Range("K17").Value = 2.222222222222
$ws.Range("L17").Value = -6.122448979591
$ws.Range("M17").Value = 39.393939393939
$ws.Range("N17").Value = -47.727272727272
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -57.142857142857
$ws.Range("I18").Value = 6
$ws.Range("J18").Value = 9
$ws.Range("K18").Value = -33.333333333333
$ws.Range("L18").Value = -40
$ws.Range("M18").Value = -87.5
$ws.Range("N18").Value = -94.690265486725
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -42.857142857142
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = 8.333333333333
$ws.Range("I19").Value = 46
$ws.Range("J19").Value = 36
$ws.Range("K19").Value = 27.777777777777
$ws.Range("L19").Value = 21.052631578947
$ws.Range("M19").Value = -41.025641025641
$ws.Range("N19").Value = -79.185520361991
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -60
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 24
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = -17.241379310344
$ws.Range("L20").Value = 20
$ws.Range("M20").Value = -29.411764705882
$ws.Range("N20").Value = -89.041095890411
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -43.478260869565
$ws.Range("F21").Value = 75
$ws.Range("G21").Value = 83
$ws.Range("H21").Value = -9.638554216867
$ws.Range("I21").Value = 136
$ws.Range("J21").Value = 137
$ws.Range("K21").Value = -0.729927007299
$ws.Range("L21").Value = -0.729927007299
$ws.Range("M21").Value = -42.127659574468
$ws.Range("N21").Value = -82.957393483709
$ws.Range("C23").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 1
$ws.Range("K23").Value = -50
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -66.666666666666
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -41.666666666666
$ws.Range("F24").Value = 58
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = -30.120481927710
$ws.Range("I24").Value = 101
$ws.Range("J24").Value = 139
$ws.Range("K24").Value = -27.338129496402
$ws.Range("L24").Value = -24.626865671641
$ws.Range("M24").Value = -12.931034482758
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -60
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = -28.571428571428
$ws.Range("I25").Value = 24
$ws.Range("J25").Value = 47
$ws.Range("K25").Value = -48.936170212766
$ws.Range("L25").Value = -36.842105263157
$ws.Range("C26").Value = 24
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 166.666666666667
$ws.Range("F26").Value = 53
$ws.Range("G26").Value = 49
$ws.Range("H26").Value = 8.163265306122
$ws.Range("I26").Value = 84
$ws.Range("J26").Value = 69
$ws.Range("K26").Value = 21.739130434782
$ws.Range("L26").Value = 16.666666666666
$ws.Range("M26").Value = -5.617977528089
$ws.Range("L27").Value = -77.777777777777
$ws.Range("C28").Value = 4
$ws.Range("F28").Value = 10
$ws.Range("H28").Value = 400
$ws.Range("I28").Value = 13
$ws.Range("K28").Value = 160
$ws.Range("L28").Value = 18.181818181818
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = -50
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -87.5
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = -50
$ws.Range("M30").Value = -40
$ws.Range("N30").Value = -84.210526315789
